# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.01253208636536152, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 0, 1.571903811587466)
    3 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 5.582307763322248)
    4 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 5.582307763322248)
    5 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 0, 8.656069925401464)
    6 = @(0.6545652718822623, 0.002658071450198252, 18.71679738969934, 0.5333859586016987, 0, 19.9074066916335)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("F$row").Value = $vals[4]
    $ws.Range("G$row").Value = $vals[5]
}
